# Update the "Cost" class-attribute label to "Price" in the Model
# Component Class Diagram (docs/diagrams/ModelComponentClassDiagram.pptx).
#
# The shape is a small "Rectangle 8" box (shape Id 83) attached to the
# SavedAmount class via an elbow connector, on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 83) {
        $target = $shp
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate shape with Id 83 (Cost label) on slide 1"
}

if ($target.HasTextFrame -and $target.TextFrame.HasText -and $target.TextFrame.TextRange.Text -eq "Cost") {
    $target.TextFrame.TextRange.Text = "Price"
}
